$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Swap rows 9 and 10 (columns B through AC) ---
$row9 = @()
$row10 = @()
for ($col = 2; $col -le 29; $col++) {
    $row9 += ,$ws.Cells.Item(9, $col).Value()
    $row10 += ,$ws.Cells.Item(10, $col).Value()
}
for ($i = 0; $i -lt $row9.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(9, $col).Value = $row10[$i]
    $ws.Cells.Item(10, $col).Value = $row9[$i]
}

# --- Step 2: Update rows 122-125 with results + refreshed odds ---
# Row 122
$ws.Cells.Item(122, 1).Value = 120
$ws.Cells.Item(122, 2).Value = 6814427
$ws.Cells.Item(122, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(122, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(122, 5).Value = 45364.44791666666
$ws.Cells.Item(122, 6).Value = "NS Mura"
$ws.Cells.Item(122, 7).Value = "NK Bravo"
$ws.Cells.Item(122, 8).Value = 1
$ws.Cells.Item(122, 9).Value = 2
$ws.Cells.Item(122, 10).Value = "A"
$ws.Cells.Item(122, 11).Value = 2.5
$ws.Cells.Item(122, 12).Value = 3.2
$ws.Cells.Item(122, 13).Value = 2.55
$ws.Cells.Item(122, 14).Value = 2.7
$ws.Cells.Item(122, 15).Value = 3
$ws.Cells.Item(122, 16).Value = 2.45
$ws.Cells.Item(122, 17).Value = 0
$ws.Cells.Item(122, 18).Value = 1.975
$ws.Cells.Item(122, 19).Value = 1.825
$ws.Cells.Item(122, 20).Value = 2
$ws.Cells.Item(122, 21).Value = 1.9
$ws.Cells.Item(122, 22).Value = 1.9
$ws.Cells.Item(122, 23).Value = -1
$ws.Cells.Item(122, 24).Value = -1
$ws.Cells.Item(122, 25).Value = 1.45
$ws.Cells.Item(122, 26).Value = -1
$ws.Cells.Item(122, 27).Value = 0.825
$ws.Cells.Item(122, 28).Value = 0.8999999999999999
$ws.Cells.Item(122, 29).Value = -1

# Row 123
$ws.Cells.Item(123, 1).Value = 121
$ws.Cells.Item(123, 2).Value = 6816448
$ws.Cells.Item(123, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(123, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(123, 5).Value = 45364.53125
$ws.Cells.Item(123, 6).Value = "NK Aluminij"
$ws.Cells.Item(123, 7).Value = "NK Rogaska"
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 9).Value = 1
$ws.Cells.Item(123, 10).Value = "A"
$ws.Cells.Item(123, 11).Value = 2.15
$ws.Cells.Item(123, 12).Value = 3.2
$ws.Cells.Item(123, 13).Value = 3
$ws.Cells.Item(123, 14).Value = 3.5
$ws.Cells.Item(123, 15).Value = 3.1
$ws.Cells.Item(123, 16).Value = 2
$ws.Cells.Item(123, 17).Value = 0.25
$ws.Cells.Item(123, 18).Value = 2.025
$ws.Cells.Item(123, 19).Value = 1.775
$ws.Cells.Item(123, 20).Value = 2.5
$ws.Cells.Item(123, 21).Value = 1.925
$ws.Cells.Item(123, 22).Value = 1.875
$ws.Cells.Item(123, 23).Value = -1
$ws.Cells.Item(123, 24).Value = -1
$ws.Cells.Item(123, 25).Value = 1
$ws.Cells.Item(123, 26).Value = -1
$ws.Cells.Item(123, 27).Value = 0.7749999999999999
$ws.Cells.Item(123, 28).Value = -1
$ws.Cells.Item(123, 29).Value = 0.875

# Row 124
$ws.Cells.Item(124, 1).Value = 122
$ws.Cells.Item(124, 2).Value = 6814426
$ws.Cells.Item(124, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(124, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(124, 5).Value = 45364.61458333334
$ws.Cells.Item(124, 6).Value = "NK Radomlje"
$ws.Cells.Item(124, 7).Value = "NK Domzale"
$ws.Cells.Item(124, 8).Value = 2
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = "H"
$ws.Cells.Item(124, 11).Value = 2.55
$ws.Cells.Item(124, 12).Value = 3.25
$ws.Cells.Item(124, 13).Value = 2.4
$ws.Cells.Item(124, 14).Value = 2.3
$ws.Cells.Item(124, 15).Value = 3.1
$ws.Cells.Item(124, 16).Value = 2.8
$ws.Cells.Item(124, 17).Value = -0.25
$ws.Cells.Item(124, 18).Value = 2.025
$ws.Cells.Item(124, 19).Value = 1.775
$ws.Cells.Item(124, 20).Value = 2.5
$ws.Cells.Item(124, 21).Value = 1.975
$ws.Cells.Item(124, 22).Value = 1.825
$ws.Cells.Item(124, 23).Value = 1.3
$ws.Cells.Item(124, 24).Value = -1
$ws.Cells.Item(124, 25).Value = -1
$ws.Cells.Item(124, 26).Value = 1.025
$ws.Cells.Item(124, 27).Value = -1
$ws.Cells.Item(124, 28).Value = -1
$ws.Cells.Item(124, 29).Value = 0.825

# Row 125
$ws.Cells.Item(125, 1).Value = 123
$ws.Cells.Item(125, 2).Value = 6814425
$ws.Cells.Item(125, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(125, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(125, 5).Value = 45365.53125
$ws.Cells.Item(125, 6).Value = "FC Koper"
$ws.Cells.Item(125, 7).Value = "NK Celje"
$ws.Cells.Item(125, 8).Value = 1
$ws.Cells.Item(125, 9).Value = 3
$ws.Cells.Item(125, 10).Value = "A"
$ws.Cells.Item(125, 11).Value = 4
$ws.Cells.Item(125, 12).Value = 3.25
$ws.Cells.Item(125, 13).Value = 1.8
$ws.Cells.Item(125, 14).Value = 4
$ws.Cells.Item(125, 15).Value = 3.4
$ws.Cells.Item(125, 16).Value = 1.75
$ws.Cells.Item(125, 17).Value = 0.5
$ws.Cells.Item(125, 18).Value = 1.975
$ws.Cells.Item(125, 19).Value = 1.825
$ws.Cells.Item(125, 20).Value = 2.5
$ws.Cells.Item(125, 21).Value = 1.85
$ws.Cells.Item(125, 22).Value = 1.95
$ws.Cells.Item(125, 23).Value = -1
$ws.Cells.Item(125, 24).Value = -1
$ws.Cells.Item(125, 25).Value = 0.75
$ws.Cells.Item(125, 26).Value = -1
$ws.Cells.Item(125, 27).Value = 0.825
$ws.Cells.Item(125, 28).Value = 0.8500000000000001
$ws.Cells.Item(125, 29).Value = -1

# --- Step 3: Append new rows 126-129 ---
$ws.Range("A125").Copy()
$ws.Range("A126:A129").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("E125").Copy()
$ws.Range("E126:E129").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 126
$ws.Cells.Item(126, 1).Value = 124
$ws.Cells.Item(126, 2).Value = 6814428
$ws.Cells.Item(126, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(126, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(126, 5).Value = 45367.45833333334
$ws.Cells.Item(126, 6).Value = "NK Bravo"
$ws.Cells.Item(126, 7).Value = "Olimpija Ljubljana"
$ws.Cells.Item(126, 8).Value = 1
$ws.Cells.Item(126, 9).Value = 1
$ws.Cells.Item(126, 10).Value = "D"
$ws.Cells.Item(126, 11).Value = 3.4
$ws.Cells.Item(126, 12).Value = 3.4
$ws.Cells.Item(126, 13).Value = 1.909
$ws.Cells.Item(126, 14).Value = 3.6
$ws.Cells.Item(126, 15).Value = 3.3
$ws.Cells.Item(126, 16).Value = 1.909
$ws.Cells.Item(126, 17).Value = 0.5
$ws.Cells.Item(126, 18).Value = 1.8
$ws.Cells.Item(126, 19).Value = 2
$ws.Cells.Item(126, 20).Value = 2.25
$ws.Cells.Item(126, 21).Value = 1.8
$ws.Cells.Item(126, 22).Value = 2
$ws.Cells.Item(126, 23).Value = -1
$ws.Cells.Item(126, 24).Value = 2.3
$ws.Cells.Item(126, 25).Value = -1
$ws.Cells.Item(126, 26).Value = 0.8
$ws.Cells.Item(126, 27).Value = -1
$ws.Cells.Item(126, 28).Value = -0.5
$ws.Cells.Item(126, 29).Value = 0.5

# Row 127
$ws.Cells.Item(127, 1).Value = 125
$ws.Cells.Item(127, 2).Value = 6814751
$ws.Cells.Item(127, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(127, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(127, 5).Value = 45367.67708333334
$ws.Cells.Item(127, 6).Value = "NK Maribor"
$ws.Cells.Item(127, 7).Value = "NK Aluminij"
$ws.Cells.Item(127, 8).Value = 2
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = "H"
$ws.Cells.Item(127, 11).Value = 1.727
$ws.Cells.Item(127, 12).Value = 3.6
$ws.Cells.Item(127, 13).Value = 3.9
$ws.Cells.Item(127, 14).Value = 1.8
$ws.Cells.Item(127, 15).Value = 3.6
$ws.Cells.Item(127, 16).Value = 3.75
$ws.Cells.Item(127, 17).Value = -0.5
$ws.Cells.Item(127, 18).Value = 1.925
$ws.Cells.Item(127, 19).Value = 1.875
$ws.Cells.Item(127, 20).Value = 2.75
$ws.Cells.Item(127, 21).Value = 1.975
$ws.Cells.Item(127, 22).Value = 1.825
$ws.Cells.Item(127, 23).Value = 0.8
$ws.Cells.Item(127, 24).Value = -1
$ws.Cells.Item(127, 25).Value = -1
$ws.Cells.Item(127, 26).Value = 0.925
$ws.Cells.Item(127, 27).Value = -1
$ws.Cells.Item(127, 28).Value = -1
$ws.Cells.Item(127, 29).Value = 0.825

# Row 128
$ws.Cells.Item(128, 1).Value = 126
$ws.Cells.Item(128, 2).Value = 6814429
$ws.Cells.Item(128, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(128, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(128, 5).Value = 45368.45833333334
$ws.Cells.Item(128, 6).Value = "NK Celje"
$ws.Cells.Item(128, 7).Value = "NK Radomlje"
$ws.Cells.Item(128, 11).Value = 1.285
$ws.Cells.Item(128, 12).Value = 5
$ws.Cells.Item(128, 13).Value = 7.5
$ws.Cells.Item(128, 14).Value = 1.363
$ws.Cells.Item(128, 15).Value = 4.5
$ws.Cells.Item(128, 16).Value = 7
$ws.Cells.Item(128, 17).Value = -1.25
$ws.Cells.Item(128, 18).Value = 1.8
$ws.Cells.Item(128, 19).Value = 2
$ws.Cells.Item(128, 20).Value = 2.75
$ws.Cells.Item(128, 21).Value = 1.95
$ws.Cells.Item(128, 22).Value = 1.85
$ws.Cells.Item(128, 23).Value = 0
$ws.Cells.Item(128, 24).Value = 0
$ws.Cells.Item(128, 25).Value = 0
$ws.Cells.Item(128, 26).Value = 0
$ws.Cells.Item(128, 27).Value = 0

# Row 129
$ws.Cells.Item(129, 1).Value = 127
$ws.Cells.Item(129, 2).Value = 6814430
$ws.Cells.Item(129, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(129, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(129, 5).Value = 45368.5625
$ws.Cells.Item(129, 6).Value = "NK Maribor"
$ws.Cells.Item(129, 7).Value = "FC Koper"
$ws.Cells.Item(129, 11).Value = 1.666
$ws.Cells.Item(129, 12).Value = 3.6
$ws.Cells.Item(129, 13).Value = 4.2
$ws.Cells.Item(129, 14).Value = 1.8
$ws.Cells.Item(129, 15).Value = 3.5
$ws.Cells.Item(129, 16).Value = 3.75
$ws.Cells.Item(129, 17).Value = -0.5
$ws.Cells.Item(129, 18).Value = 1.825
$ws.Cells.Item(129, 19).Value = 1.975
$ws.Cells.Item(129, 20).Value = 2.5
$ws.Cells.Item(129, 21).Value = 1.825
$ws.Cells.Item(129, 22).Value = 1.975
$ws.Cells.Item(129, 23).Value = 0
$ws.Cells.Item(129, 24).Value = 0
$ws.Cells.Item(129, 25).Value = 0
$ws.Cells.Item(129, 26).Value = 0
$ws.Cells.Item(129, 27).Value = 0

Write-Host "Edit complete"